$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(420).Insert()

$ws.Cells.Item(420, 1).Value = 4
$ws.Cells.Item(420, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(420, 3).Value = "Los Lagos"
$ws.Cells.Item(420, 4).Value = 45166
$ws.Cells.Item(420, 5).Value = 10
$ws.Cells.Item(420, 6).Value = 100112032
$ws.Cells.Item(420, 7).Value = "Zapallo italiano"
$ws.Cells.Item(420, 8).Value = "Sin especificar"
$ws.Cells.Item(420, 9).Value = "Primera"
$ws.Cells.Item(420, 10).Value = 70
$ws.Cells.Item(420, 11).Value = 18000
$ws.Cells.Item(420, 12).Value = 18000
$ws.Cells.Item(420, 13).Value = 18000
$ws.Cells.Item(420, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(420, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(420, 16).Value = 360
$ws.Cells.Item(420, 17).Value = 50
$ws.Cells.Item(420, 18).Value = "Hortaliza"
